$d = $word.ActiveDocument

# Paragraph 2: (Ref-DJ49KL) -> (Ref-f840991)
$r = $d.Paragraphs.Item(2).Range
$r.Find.Execute("Ref-DJ49KL", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-f840991", 2)

# Paragraph 3: (Ref-A1B2C3), (Ref-D4E5F6), (Ref-G7H8I9) -> (Ref-s317977) each
$r = $d.Paragraphs.Item(3).Range
$r.Find.Execute("Ref-A1B2C3", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s317977", 2)

$r = $d.Paragraphs.Item(3).Range
$r.Find.Execute("Ref-D4E5F6", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s317977", 2)

$r = $d.Paragraphs.Item(3).Range
$r.Find.Execute("Ref-G7H8I9", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s317977", 2)

# Paragraph 4: (Ref-A1B2C3), (Ref-D4E5F6) -> (Ref-s451092) each
$r = $d.Paragraphs.Item(4).Range
$r.Find.Execute("Ref-A1B2C3", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s451092", 2)

$r = $d.Paragraphs.Item(4).Range
$r.Find.Execute("Ref-D4E5F6", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s451092", 2)

# Paragraph 5: (Ref-DJ49F2) -> (Ref-u946651)
$r = $d.Paragraphs.Item(5).Range
$r.Find.Execute("Ref-DJ49F2", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-u946651", 2)
